$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Senders MAC" / "Senders IP" entries to the possessive form,
# and "Length" (IP total length field) to "Total length".
$ws.Range("D25").Value = "Sender's MAC"
$ws.Range("D26").Value = "Sender's IP"
$ws.Range("D12").Value = "Total length"

# Move the active selection from F22 to F5.
$ws.Range("F5").Select()
